$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated transition-probability matrix values (recomputed after simulating more games)
$ws.Range("B2").Value = 0.2442528735632184
$ws.Range("C2").Value = 0.4597701149425287
$ws.Range("J2").Value = 0.04310344827586207
$ws.Range("P2").Value = 0.1551724137931035
$ws.Range("S2").Value = 0.09770114942528736
$ws.Range("B3").Value = 0.0125
$ws.Range("C3").Value = 0.01875
$ws.Range("J3").Value = 0.06875000000000001
$ws.Range("P3").Value = 0.66875
$ws.Range("S3").Value = 0.23125
$ws.Range("J4").Value = 0.05084745762711865
$ws.Range("P4").Value = 0.711864406779661
$ws.Range("S4").Value = 0.2372881355932203
$ws.Range("B6").Value = 0.07303370786516854
$ws.Range("D6").Value = 0.01123595505617977
$ws.Range("E6").Value = 0.005617977528089887
$ws.Range("F6").Value = 0.05056179775280899
$ws.Range("J6").Value = 0.2022471910112359
$ws.Range("O6").Value = 0.03370786516853932
$ws.Range("Q6").Value = 0.1966292134831461
$ws.Range("R6").Value = 0.05056179775280899
$ws.Range("S6").Value = 0.3764044943820224
$ws.Range("B7").Value = 0.1206896551724138
$ws.Range("D7").Value = 0.02298850574712644
$ws.Range("F7").Value = 0.05747126436781609
$ws.Range("J7").Value = 0.132183908045977
$ws.Range("O7").Value = 0.005747126436781609
$ws.Range("Q7").Value = 0.1954022988505747
$ws.Range("R7").Value = 0.04022988505747126
$ws.Range("S7").Value = 0.4252873563218391
$ws.Range("B8").Value = 0.1039260969976905
$ws.Range("D8").Value = 0.03926096997690531
$ws.Range("F8").Value = 0.06466512702078522
$ws.Range("J8").Value = 0.09699769053117784
$ws.Range("O8").Value = 0.006928406466512702
$ws.Range("Q8").Value = 0.1939953810623557
$ws.Range("R8").Value = 0.1062355658198614
$ws.Range("S8").Value = 0.3879907621247113
$ws.Range("B9").Value = 0.08928571428571429
$ws.Range("D9").Value = 0.03125
$ws.Range("F9").Value = 0.05357142857142857
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.004464285714285714
$ws.Range("Q9").Value = 0.2053571428571428
$ws.Range("R9").Value = 0.1026785714285714
$ws.Range("S9").Value = 0.3883928571428572
$ws.Range("B10").Value = 0.1358234295415959
$ws.Range("D10").Value = 0.02461799660441426
$ws.Range("F10").Value = 0.06366723259762309
$ws.Range("J10").Value = 0.1129032258064516
$ws.Range("O10").Value = 0.02716468590831918
$ws.Range("Q10").Value = 0.2147707979626486
$ws.Range("R10").Value = 0.07130730050933787
$ws.Range("S10").Value = 0.3497453310696095
$ws.Range("G11").Value = 0.1602787456445993
$ws.Range("J11").Value = 0.1045296167247387
$ws.Range("K11").Value = 0.2299651567944251
$ws.Range("L11").Value = 0.4808362369337979
$ws.Range("S11").Value = 0.02439024390243903
$ws.Range("G12").Value = 0.6906474820143885
$ws.Range("J12").Value = 0.2158273381294964
$ws.Range("L12").Value = 0.01438848920863309
$ws.Range("S12").Value = 0.07913669064748201
$ws.Range("G13").Value = 0.66
$ws.Range("J13").Value = 0.3
$ws.Range("S13").Value = 0.04
$ws.Range("F15").Value = 0.005
$ws.Range("H15").Value = 0.12
$ws.Range("I15").Value = 0.075
$ws.Range("J15").Value = 0.405
$ws.Range("K15").Value = 0.05
$ws.Range("M15").Value = 0.015
$ws.Range("O15").Value = 0.045
$ws.Range("S15").Value = 0.285
$ws.Range("F16").Value = 0.01470588235294118
$ws.Range("H16").Value = 0.1715686274509804
$ws.Range("I16").Value = 0.06372549019607843
$ws.Range("J16").Value = 0.3872549019607843
$ws.Range("K16").Value = 0.1617647058823529
$ws.Range("M16").Value = 0.02450980392156863
$ws.Range("O16").Value = 0.05392156862745098
$ws.Range("S16").Value = 0.1225490196078431
$ws.Range("F17").Value = 0.0111358574610245
$ws.Range("H17").Value = 0.1737193763919822
$ws.Range("I17").Value = 0.1269487750556793
$ws.Range("J17").Value = 0.4187082405345212
$ws.Range("K17").Value = 0.0779510022271715
$ws.Range("M17").Value = 0.0155902004454343
$ws.Range("N17").Value = 0.0022271714922049
$ws.Range("O17").Value = 0.0556792873051225
$ws.Range("S17").Value = 0.1180400890868597
$ws.Range("F18").Value = 0.005882352941176471
$ws.Range("H18").Value = 0.1529411764705882
$ws.Range("I18").Value = 0.1176470588235294
$ws.Range("J18").Value = 0.3352941176470588
$ws.Range("K18").Value = 0.1470588235294118
$ws.Range("M18").Value = 0.03529411764705882
$ws.Range("O18").Value = 0.05882352941176471
$ws.Range("S18").Value = 0.1470588235294118
$ws.Range("F19").Value = 0.00919732441471572
$ws.Range("H19").Value = 0.225752508361204
$ws.Range("I19").Value = 0.1020066889632107
$ws.Range("J19").Value = 0.3561872909698997
$ws.Range("K19").Value = 0.1003344481605351
$ws.Range("M19").Value = 0.02591973244147157
$ws.Range("N19").Value = 0.0008361204013377926
$ws.Range("O19").Value = 0.06438127090301003
$ws.Range("S19").Value = 0.1153846153846154
